# Upgrade to new masterlist: fill in missing sequence values in column D
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D6").Value = 8
$ws.Range("D7").Value = 9
